$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Opening paragraph: tidy the "armed forces" wording and rewrite the
#    Coast Guard sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "does not maintain an armed forces, but relies",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "does not maintain an armed force but relies",
    2)

$d.Content.Find.Execute(
    "The only force Iceland controls is its Coast Guard.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The only military type element in Iceland is its Coast Guard.",
    2)

# ---------------------------------------------------------------------
# 2. Drop the stray empty paragraph that sits between the "consists of
#    two patrol ships..." paragraph and the Aegir/Tyr paragraph.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [char]13 -and $i -lt $d.Paragraphs.Count) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.StartsWith("The two patrol ships of the")) {
            $p.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# 3. Consolidate the runs/proofErr wraps around Aegir / Tyr by replacing
#    the text in-place (identical content forces Word to merge runs).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The two patrol ships of the " + [char]0x00C6 + [char]0x0067 + "ir class are the " + [char]0x00C6 + [char]0x0067 + "ir itself and the T" + [char]0x00FD + "r, they are only lightly armed but can support one of the four SA-365N-1 Dauphin 2 helicopters",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The two patrol ships of the " + [char]0x00C6 + [char]0x0067 + "ir class are the " + [char]0x00C6 + [char]0x0067 + "ir itself and the T" + [char]0x00FD + "r, they are only lightly armed but can support one of the four SA-365N-1 Dauphin 2 helicopters",
    2)

# Wrap the helicopter model name in curly quotes.
$d.Content.Find.Execute(
    "can support one of the four SA-365N-1 Dauphin 2 helicopters",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can support one of the four " + [char]0x2018 + "SA-365N-1 Dauphin 2" + [char]0x2019 + " helicopters",
    2)

# ---------------------------------------------------------------------
# 4. Remove the leftover _GoBack bookmark after the Tyr picture.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 5. Merge the lone bold-space paragraph into the following "The single
#    fixed wing asset..." paragraph by deleting the mark between them.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -eq (" " + [char]13) -and $i -lt $d.Paragraphs.Count) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.StartsWith("The single fixed wing asset")) {
            $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
            $markRange.Delete()
            break
        }
    }
}
